$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Match column A's font/style to column B's (s=4 -> s=3) for rows 2-7 ---
$ws.Range("B2").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: JPY=X/Macro/100/10/150/AV/1d -> ETH-USD/Indicators/150/10/300/Binance/1d ---
$ws.Range("A2").Value = "ETH-USD"
$ws.Range("B2").Value = "Indicators"
$ws.Range("C2").Value = 150
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 300
$ws.Range("F2").Value = "Binance"
$ws.Range("G2").Value = "1d"

# --- Row 3: JPY=X/Macro/150/10/150/AV/1d -> ETH-USD/Indicators/150/10/300/Binance/1d ---
$ws.Range("A3").Value = "ETH-USD"
$ws.Range("B3").Value = "Indicators"
$ws.Range("C3").Value = 150
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 300
$ws.Range("F3").Value = "Binance"
$ws.Range("G3").Value = "1d"

# --- Rows 4-7: clear out the remaining sample rows (keep formatting/styles) ---
$ws.Range("A4:G7").ClearContents()

# --- Update the selection to match the new active range ---
$null = $ws.Range("A4:G6").Select()
